# Timp3-Agtr2 LR-pairs worksheet refresh ("update scripts wuth new tpm").
#
# The upstream NATMI script was re-run with new TPM input data. For this
# ligand-receptor pair the result table now:
#   - only keeps rows where the Target cluster is "FAPs" (the old
#     Target="ECs" rows are gone),
#   - has refreshed numeric columns (E..T) for every sending cluster, and
#   - gains a brand-new sending cluster, "Inflammatory-Mac", that did not
#     exist in the previous run.
# Net effect: the data block shrinks from 8 rows (A1:T9) to 5 rows
# (A1:T6): ECs, FAPs, Inflammatory-Mac, MuSCs, Neutrophils - each paired
# with Timp3 -> Agtr2 -> FAPs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T numeric metrics.
$rows = @(
    @{ Row = 2; A = "ECs";              E = 2; F = 1;                  G = 122.1023975;         H = 244.204795;   I = 0.5529388544589152;  J = 0.5365853648354085;  K = 3; L = 1; M = 0.9663516666666667; N = 2.899055; O = 1; P = 1; Q = 117.9938553281208;  R = 707.9631319687251; S = 0.5529388544589152;  T = 0.5365853648354085 },
    @{ Row = 3; A = "FAPs";             E = 3; F = 1;                  G = 13.28375266666667;   H = 39.851258;    I = 0.06015527240095533; J = 0.08756421761939603; K = 3; L = 1; M = 0.9663516666666667; N = 2.899055; O = 1; P = 1; Q = 12.83677652902111;  R = 115.53098876119;   S = 0.06015527240095533; T = 0.08756421761939603 },
    @{ Row = 4; A = "Inflammatory-Mac"; E = 1; F = 0.3333333333333333; G = 0.02528933333333333; H = 0.075868;     I = 0.0001145223622931973; J = 0.0001667029447940724; K = 3; L = 1; M = 0.9663516666666667; N = 2.899055; O = 1; P = 1; Q = 0.02443838941555556; R = 0.21994550474;      S = 0.0001145223622931973; T = 0.0001667029447940724 },
    @{ Row = 5; A = "MuSCs";            E = 2; F = 1;                  G = 85.26190199999999;   H = 170.523804;   I = 0.3861072295559822;  J = 0.3746878826947755;  K = 3; L = 1; M = 0.9663516666666667; N = 2.899055; O = 1; P = 1; Q = 82.39298110087;     R = 494.35788660522;   S = 0.3861072295559822;  T = 0.3746878826947755 },
    @{ Row = 6; A = "Neutrophils";      E = 2; F = 0.6666666666666666; G = 0.1510706666666667;  H = 0.453212;     I = 0.0006841212218540695; J = 0.0009958319056257072; K = 3; L = 1; M = 0.9663516666666667; N = 2.899055; O = 1; P = 1; Q = 0.1459873905177778;  R = 1.31388651466;     S = 0.0006841212218540695; T = 0.0009958319056257072 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("A$i").Value = $r.A
    $ws.Range("B$i").Value = "Timp3"
    $ws.Range("C$i").Value = "Agtr2"
    $ws.Range("D$i").Value = "FAPs"
    $ws.Range("E$i").Value = $r.E
    $ws.Range("F$i").Value = $r.F
    $ws.Range("G$i").Value = $r.G
    $ws.Range("H$i").Value = $r.H
    $ws.Range("I$i").Value = $r.I
    $ws.Range("J$i").Value = $r.J
    $ws.Range("K$i").Value = $r.K
    $ws.Range("L$i").Value = $r.L
    $ws.Range("M$i").Value = $r.M
    $ws.Range("N$i").Value = $r.N
    $ws.Range("O$i").Value = $r.O
    $ws.Range("P$i").Value = $r.P
    $ws.Range("Q$i").Value = $r.Q
    $ws.Range("R$i").Value = $r.R
    $ws.Range("S$i").Value = $r.S
    $ws.Range("T$i").Value = $r.T
}

# Drop the now-obsolete rows 7-9 (old Target="ECs" data no longer present
# in the refreshed run) so the worksheet's used range / dimension shrinks
# from A1:T9 down to A1:T6, matching the new result table.
$ws.Range("A7:T9").ClearContents()
